$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Batch E1 (rows 5-12): add a new "bt2_file" value in column P ---
# hpgl0564 - hpgl0566, hpgl0567 have a trailing space difference per sample;
# reproduce the exact text (including trailing spaces) recorded for the edit.
$ws.Range("P5").Value = "preprocessing/E1/hpgl0564/outputs/bowtie2_scerevisiae/hpgl0564_forward-trimmed.count.xz"
$ws.Range("P6").Value = "preprocessing/E1/hpgl0565/outputs/bowtie2_scerevisiae/hpgl0565_forward-trimmed.count.xz"
$ws.Range("P7").Value = "preprocessing/E1/hpgl0566/outputs/bowtie2_scerevisiae/hpgl0566_forward-trimmed.count.xz"
$ws.Range("P8").Value = "preprocessing/E1/hpgl0567/outputs/bowtie2_scerevisiae/hpgl0567_forward-trimmed.count.xz "
$ws.Range("P9").Value = "preprocessing/E1/hpgl0568/outputs/bowtie2_scerevisiae/hpgl0568_forward-trimmed.count.xz "
$ws.Range("P10").Value = "preprocessing/E1/hpgl0569/outputs/bowtie2_scerevisiae/hpgl0569_forward-trimmed.count.xz "
$ws.Range("P11").Value = "preprocessing/E1/hpgl0570/outputs/bowtie2_scerevisiae/hpgl0570_forward-trimmed.count.xz "
$ws.Range("P12").Value = "preprocessing/E1/hpgl0571/outputs/bowtie2_scerevisiae/hpgl0571_forward-trimmed.count.xz "

# --- Batch E2 (rows 13-28): fix typeo "preprocessing/v2/..." -> "preprocessing/E2/..." in column P ---
$ws.Range("P13").Value = "preprocessing/E2/hpgl0774/outputs/bowtie2_scerevisiae/hpgl0774_forward-trimmed.count.xz"
$ws.Range("P14").Value = "preprocessing/E2/hpgl0775/outputs/bowtie2_scerevisiae/hpgl0775_forward-trimmed.count.xz"
$ws.Range("P15").Value = "preprocessing/E2/hpgl0776/outputs/bowtie2_scerevisiae/hpgl0776_forward-trimmed.count.xz"
$ws.Range("P16").Value = "preprocessing/E2/hpgl0777/outputs/bowtie2_scerevisiae/hpgl0777_forward-trimmed.count.xz"
$ws.Range("P17").Value = "preprocessing/E2/hpgl0778/outputs/bowtie2_scerevisiae/hpgl0778_forward-trimmed.count.xz"
$ws.Range("P18").Value = "preprocessing/E2/hpgl0779/outputs/bowtie2_scerevisiae/hpgl0779_forward-trimmed.count.xz"
$ws.Range("P19").Value = "preprocessing/E2/hpgl0780/outputs/bowtie2_scerevisiae/hpgl0780_forward-trimmed.count.xz"
$ws.Range("P20").Value = "preprocessing/E2/hpgl0781/outputs/bowtie2_scerevisiae/hpgl0781_forward-trimmed.count.xz"
$ws.Range("P21").Value = "preprocessing/E2/hpgl0782/outputs/bowtie2_scerevisiae/hpgl0782_forward-trimmed.count.xz"
$ws.Range("P22").Value = "preprocessing/E2/hpgl0783/outputs/bowtie2_scerevisiae/hpgl0783_forward-trimmed.count.xz"
$ws.Range("P23").Value = "preprocessing/E2/hpgl0784/outputs/bowtie2_scerevisiae/hpgl0784_forward-trimmed.count.xz"
$ws.Range("P24").Value = "preprocessing/E2/hpgl0785/outputs/bowtie2_scerevisiae/hpgl0785_forward-trimmed.count.xz"
$ws.Range("P25").Value = "preprocessing/E2/hpgl0786/outputs/bowtie2_scerevisiae/hpgl0786_forward-trimmed.count.xz"
$ws.Range("P26").Value = "preprocessing/E2/hpgl0787/outputs/bowtie2_scerevisiae/hpgl0787_forward-trimmed.count.xz"
$ws.Range("P27").Value = "preprocessing/E2/hpgl0788/outputs/bowtie2_scerevisiae/hpgl0788_forward-trimmed.count.xz"
$ws.Range("P28").Value = "preprocessing/E2/hpgl0789/outputs/bowtie2_scerevisiae/hpgl0789_forward-trimmed.count.xz"

# --- Update the selected cell to P12 to match the saved view state ---
$ws.Range("P12").Select()
